# Update values in the MicroPlanResponsibilities sheet
$wb = $excel.ActiveWorkbook

$wsResp = $wb.Worksheets.Item("MicroPlanResponsibilities")

$wsResp.Range("E16").Value = 810111
$wsResp.Range("F16").Value = 95.69800000000001

$wsResp.Range("E17").Value = 11992690
$wsResp.Range("F17").Value = 1292.249

$wsResp.Range("E19").Value = 1567693
$wsResp.Range("F19").Value = 180.074

$wsResp.Range("E20").Value = 359483
$wsResp.Range("F20").Value = 42.374

$wsResp.Range("E23").Value = 1890856
$wsResp.Range("F23").Value = 208.555

$wsResp.Range("E24").Value = 1299573
$wsResp.Range("F24").Value = 145.464

# Update the MicroPlanIndex sheet: rows_cleaned count
$wsIndex = $wb.Worksheets.Item("MicroPlanIndex")
$wsIndex.Range("D2").Value = 30
